$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Change 1 - Slide 3, "Black-box effect" body placeholder:
#   "For CNNs, it's even worse" -> italicize the trailing word "worse"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(5, 1)
$tr3.Characters($para3.Start + 20, 5).Font.Italic = $true

# ---------------------------------------------------------------------------
# Change 2 - Slide 4, "Rectangle 3" footnote textbox:
#   bold the words "trade-off", "interpretability" and "performance"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(3)
$tr4 = $shp4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(1, 1)
$base4 = $para4.Start
$tr4.Characters($base4 + 75, 9).Font.Bold = $true
$tr4.Characters($base4 + 93, 16).Font.Bold = $true
$tr4.Characters($base4 + 114, 11).Font.Bold = $true

# ---------------------------------------------------------------------------
# Change 3 - Slide 8, body placeholder:
#   "Every time outputs (forward)..." -> "Every time computes outputs (forward)..."
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(4, 1)
$sub8 = $tr8.Characters($para8.Start + 6, 13)
$sub8.Text = "time computes outputs "
